# Weekly crime-data refresh for 123rd Precinct CompStat report
# (week of 8/21-8/27 -> 8/28-9/3; Volume 30 Number 34 -> 35)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: bump report week number and the covered date range ---
# (Characters() edits only the run holding the changing digits, leaving the
#  surrounding "Volume 30   Number" / "Report Covering the Week ... Through ..." text untouched)
$ws.Range("A8").Characters(21, 2).Text = "35"
$ws.Range("C9").Characters(27, 9).Text = "8/28/2023"
$ws.Range("C9").Characters(47, 9).Text = "9/3/2023"

# --- Crime-stat table refresh (rows 16-29) ---

# A few cells flip from a numeric 0/blank-style entry to the sheet's placeholder
# text ("0" / "***.*") used elsewhere for "no data this period". Copying an existing
# placeholder cell onto the target reproduces both the text and its cell style exactly.
$ws.Range("C14").Copy($ws.Range("C20"))
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("E14").Copy($ws.Range("E27"))

# Row 16 (Robbery)
$ws.Range("J16").Value = 14
$ws.Range("K16").Value = 0
$ws.Range("N16").Value = -63.157894736842

# Row 17 (Fel. Assault)
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 6
$ws.Range("H17").Value = 20
$ws.Range("I17").Value = 64
$ws.Range("J17").Value = 30
$ws.Range("K17").Value = 113.333333333333
$ws.Range("L17").Value = 166.666666666667
$ws.Range("M17").Value = 120.689655172414
$ws.Range("N17").Value = -16.883116883116

# Row 18 (Burglary)
$ws.Range("M18").Value = -45.833333333333
$ws.Range("N18").Value = -80.69306930693

# Row 19 (Gr. Larceny)
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -85.714285714285
$ws.Range("F19").Value = 16
$ws.Range("G19").Value = 25
$ws.Range("H19").Value = -36
$ws.Range("I19").Value = 195
$ws.Range("J19").Value = 189
$ws.Range("K19").Value = 3.174603174603
$ws.Range("L19").Value = 82.242990654205
$ws.Range("M19").Value = 101.030927835052
$ws.Range("N19").Value = 29.139072847682

# Row 20 (G.L.A.)
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = -61.538461538461
$ws.Range("J20").Value = 74
$ws.Range("K20").Value = -22.972972972973
$ws.Range("L20").Value = 103.571428571429
$ws.Range("M20").Value = 119.230769230769
$ws.Range("N20").Value = -88.223140495867

# Row 21 (TOTAL)
$ws.Range("C21").Value = 2
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = -85.714285714285
$ws.Range("F21").Value = 29
$ws.Range("G21").Value = 48
$ws.Range("H21").Value = -39.583333333333
$ws.Range("I21").Value = 371
$ws.Range("J21").Value = 326
$ws.Range("K21").Value = 13.803680981595
$ws.Range("L21").Value = 100.540540540541
$ws.Range("M21").Value = 53.941908713692
$ws.Range("N21").Value = -61.192468619246

# Row 24 (Petit Larceny)
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 8
$ws.Range("E24").Value = -12.5
$ws.Range("F24").Value = 32
$ws.Range("G24").Value = 51
$ws.Range("H24").Value = -37.254901960784
$ws.Range("I24").Value = 318
$ws.Range("J24").Value = 298
$ws.Range("K24").Value = 6.711409395973
$ws.Range("L24").Value = 82.758620689655
$ws.Range("M24").Value = -15.649867374005

# Row 25 (Misd. Assault)
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 36.363636363636
$ws.Range("J25").Value = 111
$ws.Range("K25").Value = 23.423423423423
$ws.Range("L25").Value = 48.91304347826
$ws.Range("M25").Value = -6.802721088435

# Row 27 (Other Sex Crimes)
$ws.Range("G27").Value = 2
$ws.Range("L27").Value = 33.333333333333

# Row 28 (Shooting Vic.)
$ws.Range("N28").Value = -66.666666666666

# Row 29 (Shooting Inc.)
$ws.Range("N29").Value = -66.666666666666

